$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --------------------------------------------------------------------------
# Header row (row 1): the sheet grows from 3 columns (QR_Data/Fecha/Hora)
# to 10 columns. Re-title the first three and add seven more headers.
# --------------------------------------------------------------------------
$ws.Range("A1").Value = "Unnamed: 0"
$ws.Range("B1").Value = "Unnamed: 1"
$ws.Range("C1").Value = "Unnamed: 2"
$ws.Range("D1").Value = "Unnamed: 3"
$ws.Range("E1").Value = "Unnamed: 4"
$ws.Range("F1").Value = "DNI"
$ws.Range("G1").Value = "Apellido"
$ws.Range("H1").Value = "Nombre"
$ws.Range("I1").Value = "Fecha"
$ws.Range("J1").Value = "Hora"

# Give the new header cells (D1:J1) the same bold/centered/bordered look
# that A1:C1 already carry.
$ws.Range("A1").Copy()
$ws.Range("D1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --------------------------------------------------------------------------
# Drop the old sample QR rows (2-5) - the sheet is being reloaded with
# fresh attendance data instead.
# --------------------------------------------------------------------------
$ws.Range("A2:C5").ClearContents()

# Helper: write a string into a cell while keeping it plain text (no
# number/date auto-conversion) and without leaving any NumberFormat/style
# residue behind - build it as a formula returning the literal text, then
# collapse the formula down to its value in place.
function Set-TextValue($addr, $text) {
    $ws.Range($addr).Formula = '="' + $text + '"'
    $ws.Range($addr).Copy()
    $ws.Range($addr).PasteSpecial(-4163)
    $excel.CutCopyMode = 0
}

# --------------------------------------------------------------------------
# Row 11 - first attendee, entered under the original A:E columns.
# --------------------------------------------------------------------------
$ws.Range("A11").Value = 44189151
Set-TextValue "B11" "Majolli"
Set-TextValue "C11" "Facundo"
Set-TextValue "D11" "2023-11-13"
Set-TextValue "E11" "02:48:43"

# --------------------------------------------------------------------------
# Row 12 - same attendee scanning again, now landing in the new DNI..Hora
# columns (F:J).
# --------------------------------------------------------------------------
$ws.Range("F12").Value = 44189151
Set-TextValue "G12" "Majolli"
Set-TextValue "H12" "Facundo"
Set-TextValue "I12" "2023-11-13"
Set-TextValue "J12" "02:55:17"

# --------------------------------------------------------------------------
# Row 13 - third scan.
# --------------------------------------------------------------------------
$ws.Range("F13").Value = 44189151
Set-TextValue "G13" "Majolli"
Set-TextValue "H13" "Facundo"
Set-TextValue "I13" "2023-11-13"
Set-TextValue "J13" "02:57:05"

# --------------------------------------------------------------------------
# Row 14 - fourth scan; this time the DNI was typed in as text, not a
# number (the data-entry slip the commit message is fixing elsewhere).
# --------------------------------------------------------------------------
Set-TextValue "F14" "44189151"
Set-TextValue "G14" "Majolli"
Set-TextValue "H14" "Facundo"
Set-TextValue "I14" "2023-11-13"
Set-TextValue "J14" "02:58:06"
